$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts every existing column
# (A->B, B->C, ..., H->I) one position to the right, carrying along values,
# formulas (re-pointed to the new references) and the merged-cell ranges.
$ws.Columns("A").Insert()

# New narrow label column (A) identifying which rows are "d" (diameter)
# router-group rows: H3, H2, H1 (rows 5, 6, 7 after the shift).
$ws.Range("A5").Value = "d"
$ws.Range("A6").Value = "d"
$ws.Range("A7").Value = "d"

# Bug fix noted in the commit message: H4's y-coordinate (raw inch value,
# now in E4 after the column insert) was left positive by mistake - it
# should be negative, matching the other router-group points.
$ws.Range("E4").Value = -1080

# Hide the raw "inch" input columns (now C and E) so only the converted
# "meter" columns (D and F) are visible, same as the author's layout.
$ws.Columns("C").Hidden = $true
$ws.Columns("E").Hidden = $true

# Move the active selection to F9 (was F8 before the column insert shifted
# the grid), matching the saved view state of the edited workbook.
$ws.Range("F9").Select()
